$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '24.427.10'
$ws.Range('E2').Value = '  -3.50%  '
$ws.Range('D3').Value = '1.644.04'
$ws.Range('E3').Value = '  -5.87%  '
Set-TextValue $ws.Range('D4') '1.001'
$ws.Range('E4').Value = '  -0.87%  '
$ws.Range('B5').Value = 'USDC'
$ws.Range('C5').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue $ws.Range('D5') '0.9991'
$ws.Range('E5').Value = '  -0.61%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue $ws.Range('D6') '305.88'
$ws.Range('E6').Value = '  -3.14%  '
Set-TextValue $ws.Range('D7') '0.3623'
$ws.Range('E7').Value = '  -4.85%  '
Set-TextValue $ws.Range('D8') '47.24'
$ws.Range('E8').Value = '  -4.44%  '
Set-TextValue $ws.Range('D9') '0.3259'
$ws.Range('E9').Value = '  -8.55%  '
Set-TextValue $ws.Range('D10') '1.113'
$ws.Range('E10').Value = '  -8.23%  '
Set-TextValue $ws.Range('D11') '0.06911'
$ws.Range('E11').Value = '  -9.28%  '
Set-TextValue $ws.Range('D12') '0.9999'
$ws.Range('E12').Value = '  -0.62%  '
Set-TextValue $ws.Range('D13') '5.938'
$ws.Range('E13').Value = '  -7.68%  '
Set-TextValue $ws.Range('D14') '19.06'
$ws.Range('E14').Value = '  -10.48%  '
$ws.Range('D15').Value = '1.644.28'
$ws.Range('E15').Value = '  -6.22%  '
$ws.Range('E16').Value = '  -7.88%  '
$ws.Range('E17').Value = '  -9.14%  '
Set-TextValue $ws.Range('D18') '0.06489'
$ws.Range('E18').Value = '  -3.82%  '
Set-TextValue $ws.Range('D19') '0.9988'
$ws.Range('E19').Value = '  -0.58%  '
Set-TextValue $ws.Range('D20') '76.66'
$ws.Range('E20').Value = '  -10.60%  '
Set-TextValue $ws.Range('D21') '5.870'
$ws.Range('E21').Value = '  -9.41%  '
Set-TextValue $ws.Range('D22') '15.67'
$ws.Range('E22').Value = '  -10.41%  '
Set-TextValue $ws.Range('D23') '12.08'
$ws.Range('E23').Value = '  -7.15%  '
$ws.Range('D24').Value = '24.377.36'
$ws.Range('E24').Value = '  -3.74%  '
Set-TextValue $ws.Range('D25') '2.395'
$ws.Range('E25').Value = '  -2.53%  '
Set-TextValue $ws.Range('D26') '2.323'
$ws.Range('E26').Value = '  -19.10%  '
Set-TextValue $ws.Range('D27') '145.31'
$ws.Range('E27').Value = '  -5.56%  '
Set-TextValue $ws.Range('D28') '18.47'
$ws.Range('E28').Value = '  -11.21%  '
$ws.Range('D29').Value = '1.828.01'
$ws.Range('E29').Value = '  -6.23%  '
Set-TextValue $ws.Range('D30') '124.55'
$ws.Range('E30').Value = '  -6.48%  '
Set-TextValue $ws.Range('D31') '1.143'
$ws.Range('E31').Value = '  -4.89%  '
Set-TextValue $ws.Range('D32') '4.050'
$ws.Range('E32').Value = '  -4.35%  '
Set-TextValue $ws.Range('D33') '5.640'
$ws.Range('E33').Value = '  -20.50%  '
Set-TextValue $ws.Range('D34') '1.681'
$ws.Range('E34').Value = '  -6.85%  '
Set-TextValue $ws.Range('D35') '0.08308'
$ws.Range('E35').Value = '  -4.66%  '
Set-TextValue $ws.Range('D36') '12.29'
$ws.Range('E36').Value = '  -13.81%  '
Set-TextValue $ws.Range('D37') '5.142'
$ws.Range('E37').Value = '  -9.95%  '
Set-TextValue $ws.Range('D38') '0.06034'
$ws.Range('E38').Value = '  -10.93%  '
Set-TextValue $ws.Range('D39') '0.02209'
$ws.Range('E39').Value = '  -10.63%  '
Set-TextValue $ws.Range('D40') '8.181'
$ws.Range('E40').Value = '  -12.25%  '
Set-TextValue $ws.Range('D41') '1.197'
$ws.Range('E41').Value = '  -6.75%  '
Set-TextValue $ws.Range('D42') '0.2033'
$ws.Range('E42').Value = '  -9.23%  '
Set-TextValue $ws.Range('D43') '0.9987'
$ws.Range('E43').Value = '  -0.53%  '
Set-TextValue $ws.Range('D44') '0.5829'
$ws.Range('E44').Value = '  -10.78%  '
Set-TextValue $ws.Range('D45') '3.716'
$ws.Range('E45').Value = '  -4.48%  '
Set-TextValue $ws.Range('D46') '12.60'
$ws.Range('E46').Value = '  -12.40%  '
Set-TextValue $ws.Range('D47') '0.5586'
$ws.Range('E47').Value = '  -10.50%  '
Set-TextValue $ws.Range('D48') '121.62'
$ws.Range('E48').Value = '  -7.19%  '
Set-TextValue $ws.Range('D49') '1.930'
$ws.Range('E49').Value = '  -11.11%  '
Set-TextValue $ws.Range('D50') '0.06887'
$ws.Range('E50').Value = '  -6.86%  '
Set-TextValue $ws.Range('D51') '73.55'
$ws.Range('E51').Value = '  -8.62%  '
